$d = $word.ActiveDocument
$p1 = $d.Paragraphs(1)
$r = $p1.Range
$r2 = $d.Range($r.Start, $r.End - 1)
$xml1 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">JSON is everything that we thought XML was going to be. FIGHT ME!! Back around 2005 when pterodactyls were still delivering the mail, this thing called XML hit the scene. Extensible Markup Language. </w:t></w:r><w:r><w:t>Sadly,</w:t></w:r><w:r><w:t xml:space="preserve"> XML failed to deliver on its promises. </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$r2.InsertXML($xml1)

$endOfP1 = $d.Paragraphs(1).Range.End
$insertPoint = $d.Range($endOfP1, $endOfP1)
$xml2 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p/><w:p><w:r><w:t>Fortunately</w:t></w:r><w:r><w:t xml:space="preserve">, along came this thing called JavaScript Object Notation or JSON. JSON is a flexible file format. The key difference between JSON and CSVs is that JSON records can be of any </w:t></w:r><w:r><w:t>length</w:t></w:r><w:r><w:t xml:space="preserve"> unlike CSVs where every row has the same number of columns. JSON records can also be nested. </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>All of</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> that means is you can&apos;t work with JSON and CSVs the same way without some transformation work. </w:t></w:r><w:r><w:t>Fortunately</w:t></w:r><w:r><w:t>, Pandas gives us those tools.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$insertPoint.InsertXML($xml2)
